# Update the workbook to match target revision.

$wb = $excel.ActiveWorkbook

# --- 1. Processes sheet: swap columns B (type) and E (description) data ---
$ws = $wb.Worksheets.Item("Processes")

$rows = 1,2,3,4,5
foreach ($r in $rows) {
    $bAddr = "B$r"
    $eAddr = "E$r"
    $bVal = $ws.Range($bAddr).Text
    $eVal = $ws.Range($eAddr).Text
    $ws.Range($bAddr).Value = $eVal
    $ws.Range($eAddr).Value = $bVal
}

# --- 2. Column widths on Processes (B narrower, E wider) ---
$ws.Columns.Item(2).ColumnWidth = 14.9
$ws.Columns.Item(5).ColumnWidth = 29.5

# --- 3. Move the process-type data validation list from E2:E5 to B2:B5 ---
$ws.Range("E2:E5").Validation.Delete()
$ws.Range("B2:B5").Validation.Add(3, 1, 1, "=Validate!`$B`$2:`$B`$3")
$ws.Range("B2:B5").Validation.ErrorTitle = "Process Type"
$ws.Range("B2:B5").Validation.ErrorMessage = "Invalid Process Type"

# --- 4. Update the cgam_processes defined name range (E4 -> D4) ---
$nm = $wb.Names.Item("cgam_processes")
$nm.RefersTo = "=Processes!`$A`$1:`$D`$4"

# --- 5. Active sheet / selection changes: Flows loses tabSelected, Processes gains it ---
$ws.Activate() | Out-Null
$ws.Range("D4").Select() | Out-Null

